$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-17 Saturday", 2)

# Update the division problems in the table.
# Row/column (1-indexed) -> new text. We resolve a fresh document Range
# from each cell's Start/End positions (rather than reusing Cell.Range
# directly) before calling Find.Execute, since duplicate cell contents
# (e.g. the two "143÷8=" cells) otherwise cause the wrong cell to be
# updated.
$t = $d.Tables(1)

function Replace-CellText($table, $row, $col, $old, $new) {
    $cell = $table.Cell($row, $col)
    $start = $cell.Range.Start
    $end = $cell.Range.End
    $r = $d.Range($start, $end)
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

$replacements = @(
    @{Row=1;  Col=1; Old="361÷5="; New="177÷6="},
    @{Row=1;  Col=2; Old="711÷2="; New="753÷9="},
    @{Row=1;  Col=3; Old="357÷2="; New="239÷8="},
    @{Row=1;  Col=4; Old="143÷8="; New="616÷7="},
    @{Row=1;  Col=5; Old="175÷4="; New="806÷9="},

    @{Row=5;  Col=1; Old="143÷8="; New="925÷6="},
    @{Row=5;  Col=2; Old="281÷6="; New="810÷7="},
    @{Row=5;  Col=3; Old="921÷3="; New="168÷2="},
    @{Row=5;  Col=4; Old="304÷6="; New="185÷8="},
    @{Row=5;  Col=5; Old="810÷8="; New="448÷7="},

    @{Row=9;  Col=1; Old="473÷2="; New="433÷4="},
    @{Row=9;  Col=2; Old="322÷8="; New="748÷9="},
    @{Row=9;  Col=3; Old="826÷7="; New="313÷4="},
    @{Row=9;  Col=4; Old="790÷7="; New="815÷3="},
    @{Row=9;  Col=5; Old="907÷5="; New="801÷9="},

    @{Row=13; Col=1; Old="172÷8="; New="194÷6="},
    @{Row=13; Col=2; Old="256÷8="; New="918÷5="},
    @{Row=13; Col=3; Old="263÷5="; New="759÷6="},
    @{Row=13; Col=4; Old="264÷9="; New="483÷9="},
    @{Row=13; Col=5; Old="605÷8="; New="808÷8="},

    @{Row=17; Col=1; Old="973÷2="; New="691÷3="},
    @{Row=17; Col=2; Old="497÷4="; New="946÷3="},
    @{Row=17; Col=3; Old="139÷8="; New="140÷8="},
    @{Row=17; Col=4; Old="671÷5="; New="194÷4="},
    @{Row=17; Col=5; Old="357÷3="; New="726÷6="}
)

foreach ($r in $replacements) {
    Replace-CellText $t $r.Row $r.Col $r.Old $r.New
}
